$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Income Statement section
$ws.Range("D8").Value = 200
$ws.Range("D10").Value = 100
$ws.Range("D17").Value = 1000
$ws.Range("D18").Value = -800
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = -800
$ws.Range("F21").Value = "NA"
$ws.Range("G21").Value = "NA"
$ws.Range("D26").Value = -800
$ws.Range("D27").Value = -800
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = -800
$ws.Range("D35").Value = -800

# Balance Sheet section
$ws.Range("D48").Value = 200
$ws.Range("D57").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("D76").Value = 900

# Cash Flow Statement section
$ws.Range("D81").Value = -800
$ws.Range("E83").Value = 0
